# VISTA-5536: Update Amenity type identifiers in the building-amenities
# lookup sheet:
#   - "furnishedUnitsAvailable"  -> "furnishedUnits"
#   - "intrusionAlarmAvailable"  -> "intrusionAlarm"
# Also update the sheet's view/selection state (scroll position & active
# cell) to reflect where the author was last working in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13 holds the "Furnished Units Available" amenity; column A stores the
# internal/enum-like type key used by the application.
$ws.Range("A13").Value = "furnishedUnits"

# Row 14 holds the "Intrusion Alarm Available" amenity; column A stores the
# internal/enum-like type key used by the application.
$ws.Range("A14").Value = "intrusionAlarm"

# Move the view/selection: scroll so row 7 is at the top-left and select A14
# (was topLeftCell A4 / selection B18 before the edit).
$ws.Range("A14").Select()
